# feat: Create add-image page
#
# 1) "Create image UI"                -> "Add" + " image UI"        (2 runs)
# 2) "Creat" / "e Image" (Strong rPr) -> "Add" / " Image" (Strong)  (text only)
# 3) "Number of characters" (numId 28 under "Create Image (admin-only)")
#                                      -> same run + new " to find" run

$d = $word.ActiveDocument

function Get-ParaIndexByExactText($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -eq ($text + "`r")) {
            return $i
        }
    }
    return -1
}

$PKG_OPEN = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$PKG_CLOSE = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Change 1: "Create image UI" -> two runs: "Add" + " image UI"
# ---------------------------------------------------------------------------
$idx1 = Get-ParaIndexByExactText("Create image UI")
if ($idx1 -eq -1) { throw "Could not locate paragraph 'Create image UI'" }
$p1 = $d.Paragraphs.Item($idx1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$xml1 = $PKG_OPEN + '<w:r><w:t>Add</w:t></w:r><w:r><w:t xml:space="preserve"> image UI</w:t></w:r>' + $PKG_CLOSE
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: "Create Image (admin-only)" heading -> "Add Image (admin-only)"
#           First run "Creat" -> "Add"; second run "e Image" -> " Image";
#           both keep their original rsid/rPr (Strong style), third run
#           " (admin-only)" is untouched.
# ---------------------------------------------------------------------------
$idx2 = Get-ParaIndexByExactText("Create Image (admin-only)")
if ($idx2 -eq -1) { throw "Could not locate paragraph 'Create Image (admin-only)'" }
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xml2 = $PKG_OPEN `
    + '<w:r w:rsidRPr="006B50A1"><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t>Add</w:t></w:r>' `
    + '<w:r w:rsidR="006137D6" w:rsidRPr="006B50A1"><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t xml:space="preserve"> Image</w:t></w:r>' `
    + '<w:r w:rsidRPr="006B50A1"><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t xml:space="preserve"> (admin-only)</w:t></w:r>' `
    + $PKG_CLOSE
$r2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: "Number of characters" (the one under the "Create Image
#           (admin-only)" field list, immediately followed by the
#           "Characters information: Name and Position (center-x, ..." para)
#           gets a new run " to find" appended after it.
# ---------------------------------------------------------------------------
$idx3 = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Number of characters`r") {
        $nxt = $p.Next()
        if ($nxt.Range.Text -like "Characters information: Name and Position (center-x*") {
            $idx3 = $i
            break
        }
    }
}
if ($idx3 -eq -1) { throw "Could not locate paragraph 'Number of characters'" }
$p3 = $d.Paragraphs.Item($idx3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End)
$xml3 = $PKG_OPEN + '<w:r><w:t>Number of characters</w:t></w:r><w:r><w:t xml:space="preserve"> to find</w:t></w:r>' + $PKG_CLOSE
$r3.InsertXML($xml3)

Write-Output "Done: para1=$idx1 para2=$idx2 para3=$idx3"
